$wb = $excel.ActiveWorkbook

# ALC!row111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 64668.625
$ws.Range("I111").Value = 2329.7
$ws.Range("J111").Value = 168566.83
$ws.Range("K111").Value = 6989.099999999999
$ws.Range("L111").Value = 505700.49
$ws.Range("M111").Value = -3922.099999999999
$ws.Range("N111").Value = -511834.49

# ALC!row131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2069.5454
$ws.Range("I131").Value = 1740
$ws.Range("J131").Value = 3552.5
$ws.Range("K131").Value = 5220
$ws.Range("L131").Value = 10657.5
$ws.Range("M131").Value = -180
$ws.Range("N131").Value = -20737.5

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2238.0977
$ws.Range("I138").Value = 775.7143
$ws.Range("J138").Value = 4409.515
$ws.Range("K138").Value = 2327.1429
$ws.Range("L138").Value = 13228.545
$ws.Range("M138").Value = 2812.8571
$ws.Range("N138").Value = -23508.545

# ALC!row139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 42580
$ws.Range("J139").Value = 42580
$ws.Range("L139").Value = 42580
$ws.Range("N139").Value = -52860

# ALC!row140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360

# ARM!row2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6759.892
$ws.Range("I32").Value = 5297.265
$ws.Range("K32").Value = 5297.265
$ws.Range("M32").Value = -5010.265

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9451.846
$ws.Range("I45").Value = 15326.286
$ws.Range("J45").Value = 2598.3333
$ws.Range("K45").Value = 15326.286
$ws.Range("L45").Value = 2598.3333
$ws.Range("M45").Value = -14949.286
$ws.Range("N45").Value = -3352.3333

# ARM!row109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 44000.2
$ws.Range("J109").Value = 44000.2
$ws.Range("L109").Value = 44000.2
$ws.Range("N109").Value = -46774.2

# ARM!row116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

# ARM!row121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# BSM!row3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

# BSM!row20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 16110.941
$ws.Range("I20").Value = 1488
$ws.Range("J20").Value = 26347
$ws.Range("K20").Value = 1488
$ws.Range("L20").Value = 26347
$ws.Range("M20").Value = -1241
$ws.Range("N20").Value = -26841

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 26318590
$ws.Range("I105").Value = 41668176
$ws.Range("J105").Value = 5015.4287
$ws.Range("K105").Value = 41668176
$ws.Range("L105").Value = 5015.4287
$ws.Range("M105").Value = -41666429
$ws.Range("N105").Value = -8509.4287

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3935.1177
$ws.Range("I134").Value = 4547.8057
$ws.Range("K134").Value = 13643.4171
$ws.Range("M134").Value = -11108.4171

# CRP!row94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3611.862
$ws.Range("I94").Value = 4315.8184
$ws.Range("J94").Value = 3181.6667
$ws.Range("K94").Value = 4315.8184
$ws.Range("L94").Value = 3181.6667
$ws.Range("M94").Value = -3864.8184
$ws.Range("N94").Value = -4083.6667

# CRP!row96
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 18177.6
$ws.Range("J96").Value = 18177.6
$ws.Range("L96").Value = 18177.6
$ws.Range("N96").Value = -23669.6

# CUL!row59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 3650
$ws.Range("I59").Value = 550
$ws.Range("J59").Value = 5200
$ws.Range("K59").Value = 1650
$ws.Range("L59").Value = 15600
$ws.Range("M59").Value = -1110
$ws.Range("N59").Value = -16680

# CUL!row60
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 385.44446
$ws.Range("I60").Value = 244.83333
$ws.Range("J60").Value = 666.6667
$ws.Range("K60").Value = 734.49999
$ws.Range("L60").Value = 2000.0001
$ws.Range("M60").Value = -483.49999
$ws.Range("N60").Value = -2502.0001

# CUL!row62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 6675.4546
$ws.Range("J62").Value = 6675.4546
$ws.Range("L62").Value = 20026.3638
$ws.Range("N62").Value = -21398.3638

# CUL!row65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 6675.4546
$ws.Range("J65").Value = 6675.4546
$ws.Range("L65").Value = 60079.0914
$ws.Range("N65").Value = -66943.0914

# CUL!row70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2302.0715
$ws.Range("I70").Value = 519.1667
$ws.Range("J70").Value = 3639.25
$ws.Range("K70").Value = 1557.5001
$ws.Range("L70").Value = 10917.75
$ws.Range("M70").Value = -1242.5001
$ws.Range("N70").Value = -11547.75

# CUL!row73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2302.0715
$ws.Range("I73").Value = 519.1667
$ws.Range("J73").Value = 3639.25
$ws.Range("K73").Value = 1557.5001
$ws.Range("L73").Value = 10917.75
$ws.Range("M73").Value = -465.5001
$ws.Range("N73").Value = -13101.75

# CUL!row75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 23810594

# CUL!row78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 23810594

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2609264.2
$ws.Range("I113").Value = 3572022
$ws.Range("J113").Value = 1111641.1
$ws.Range("K113").Value = 10716066
$ws.Range("L113").Value = 3334923.3
$ws.Range("M113").Value = -10713896
$ws.Range("N113").Value = -3339263.3

# CUL!row117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 19619346
$ws.Range("I117").Value = 17004.834
$ws.Range("K117").Value = 51014.50199999999
$ws.Range("M117").Value = -47572.50199999999

# GSM!row107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1646.7727
$ws.Range("I107").Value = 663.61536
$ws.Range("J107").Value = 3066.889
$ws.Range("K107").Value = 663.61536
$ws.Range("L107").Value = 3066.889
$ws.Range("M107").Value = 1256.38464
$ws.Range("N107").Value = -6906.889

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 125001730
$ws.Range("I113").Value = 200001200
$ws.Range("J113").Value = 2600
$ws.Range("K113").Value = 200001200
$ws.Range("L113").Value = 2600
$ws.Range("M113").Value = -199999030
$ws.Range("N113").Value = -6940

# LTW!row61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2798
$ws.Range("N61").ClearContents()

# LTW!row109
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 49750.25
$ws.Range("J109").Value = 49750.25
$ws.Range("L109").Value = 49750.25
$ws.Range("N109").Value = -52524.25

# LTW!row113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -830
$ws.Range("N113").ClearContents()

# LTW!row134
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 45415.832
$ws.Range("J134").Value = 45415.832
$ws.Range("L134").Value = 45415.832
$ws.Range("N134").Value = -55555.832

# WVR!row95
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 32660
$ws.Range("J95").Value = 32660
$ws.Range("L95").Value = 32660
$ws.Range("N95").Value = -38152

# WVR!row97
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 41999.5
$ws.Range("J97").Value = 41999.5
$ws.Range("L97").Value = 41999.5
$ws.Range("N97").Value = -43981.5

# WVR!row107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 52000384
$ws.Range("I107").Value = 83333710
$ws.Range("J107").Value = 5000396.5
$ws.Range("K107").Value = 250001130
$ws.Range("L107").Value = 15001189.5
$ws.Range("M107").Value = -249999210
$ws.Range("N107").Value = -15005029.5

# WVR!row108
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 42933.332
$ws.Range("J108").Value = 42933.332
$ws.Range("L108").Value = 42933.332
$ws.Range("N108").Value = -50613.332

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1280.0883
$ws.Range("I132").Value = 973.8261
$ws.Range("J132").Value = 1920.4546
$ws.Range("K132").Value = 2921.4783
$ws.Range("L132").Value = 5761.3638
$ws.Range("M132").Value = -391.4782999999998
$ws.Range("N132").Value = -10821.3638
